# This workbook is a weekly price log for "Plátano" (banana) at the Terminal
# La Palmera de La Serena market. Each week's data occupies three rows
# (Pintón / Primera Maduro / Primera Pintón). A new week of data (dated
# 44438, i.e. 2021-09-09) is being inserted at the top of this block
# (rows 205-207), which pushes every existing week's price row down by
# three rows (one week). The oldest existing week (previously rows
# 310-312, dated 44432) is appended as brand-new rows 313-315 at the
# bottom, preserving the sheet's chronological "sliding window" layout.
#
# Only columns D (Fecha), M (Volumen), N/O/P (Precio min/max/promedio) and
# S (Precio $/Kg) vary between weeks for a given quality row; all the
# other columns (market/product metadata) are identical down the whole
# column, so they do not need to move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Step 1: create the 3 brand-new rows (313-315) at the bottom of the
# block by duplicating rows 310-312 (values + the date format on column
# D). Their price data gets overwritten with the correct (shifted)
# values in step 2 below, together with every other existing row.
# ---------------------------------------------------------------------
for ($i = 0; $i -lt 3; $i++) {
    $destRow = 313 + $i
    $srcRow  = 310 + $i
    for ($col = 1; $col -le 20; $col++) {
        $ws.Cells.Item($destRow, $col).Value2 = $ws.Cells.Item($srcRow, $col).Value2
    }
    $ws.Cells.Item($destRow, 4).NumberFormat = $ws.Cells.Item($srcRow, 4).NumberFormat
}

# ---------------------------------------------------------------------
# Step 2: shift every existing week's price data down by one week (3
# rows). Walk from the bottom (row 315) up to row 208 so that each
# source row (r-3) is always read before it is itself overwritten.
# ---------------------------------------------------------------------
for ($r = 315; $r -ge 208; $r--) {
    $src = $r - 3
    $ws.Cells.Item($r, 4).Value2  = $ws.Cells.Item($src, 4).Value2   # Fecha
    $ws.Cells.Item($r, 13).Value2 = $ws.Cells.Item($src, 13).Value2  # Volumen
    $ws.Cells.Item($r, 14).Value2 = $ws.Cells.Item($src, 14).Value2  # Precio minimo
    $ws.Cells.Item($r, 15).Value2 = $ws.Cells.Item($src, 15).Value2  # Precio maximo
    $ws.Cells.Item($r, 16).Value2 = $ws.Cells.Item($src, 16).Value2  # Precio promedio ponderado
    $ws.Cells.Item($r, 19).Value2 = $ws.Cells.Item($src, 19).Value2  # Precio $/Kg
}

# ---------------------------------------------------------------------
# Step 3: fill in the brand-new top week (rows 205-207), dated 44438.
# ---------------------------------------------------------------------
$ws.Cells.Item(205, 4).Value2 = 44438

$ws.Cells.Item(206, 4).Value2  = 44438
$ws.Cells.Item(206, 14).Value2 = 16500
$ws.Cells.Item(206, 15).Value2 = 16500
$ws.Cells.Item(206, 16).Value2 = 16500
$ws.Cells.Item(206, 19).Value2 = 825

$ws.Cells.Item(207, 4).Value2  = 44438
$ws.Cells.Item(207, 14).Value2 = 17000
$ws.Cells.Item(207, 15).Value2 = 17000
$ws.Cells.Item(207, 16).Value2 = 17000
$ws.Cells.Item(207, 19).Value2 = 850
